$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2023-03-08 12:57:07"

# Update ratingAmount (column D) for rows whose crawl re-counted one more rating.
$ratingAmountUpdates = @{
    23  = 26
    46  = 33
    69  = 23
    88  = 18
    92  = 17
    115 = 31
    144 = 33
    204 = 13
    255 = 13
}

foreach ($row in $ratingAmountUpdates.Keys) {
    $ws.Range("D$row").Value = $ratingAmountUpdates[$row]
}

# Row 335 previously had no rating data at all; it now has a single 5-star rating.
$ws.Range("D335").Value = 1
$ws.Range("E335").Value = 5

# Every row's crawl timestamp (column O) advances to the new crawl run time.
$lastRow = 393
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("O$r").Value = $newTimestamp
}
